$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "String / Stack / Queue" filler block (rows 93-95, cols B:E)
# further down the sheet to rows 110-112 to make room for the newly solved
# Linked List problems. Doing this first so the vacated rows 93/94 are free
# to receive the new problem rows below (matches the destination row shells
# left behind by a cut, i.e. no B:E content left in rows 93-95).
$ws.Range("B93:E95").Cut($ws.Range("B110:E112"))
$ws.Range("B93:E95").Clear()

# --- Row 91: Reverse Linked List (ID 90 already present in A91)
$ws.Range("B91").Value = "Linked List"
$ws.Range("C91").Value = "Reverse Linked List"
$ws.Range("D91").Value = "Easy"
$ws.Range("E91").Value = "Done"
$ws.Range("F90").Copy()
$ws.Range("F91").PasteSpecial(-4122)
$ws.Range("F91").Value = 45922
$ws.Range("G91").Value = "O(n)"
$ws.Range("H91").Value = "O(1)"
$ws.Range("I91").Value = "Swapping next pointer"

# --- Row 92 (new): Reverse Linked List II
$ws.Range("A92").Value = 91
$ws.Range("B92").Value = "Linked List"
$ws.Range("C92").Value = "Reverse Linked List II"
$ws.Range("D92").Value = "Medium"
$ws.Range("E92").Value = "Done"
$ws.Range("F90").Copy()
$ws.Range("F92").PasteSpecial(-4122)
$ws.Range("F92").Value = 45922
$ws.Range("G92").Value = "O(n)"
$ws.Range("H92").Value = "O(1)"
$ws.Range("I92").Value = "Swapping next pointer"

# --- Row 93: Odd & Even Linked List
$ws.Range("A93").Value = 92
$ws.Range("B93").Value = "Linked List"
$ws.Range("C93").Value = "Odd & Even Linked List"
$ws.Range("D93").Value = "Medium"
$ws.Range("E93").Value = "Done"
$ws.Range("F90").Copy()
$ws.Range("F93").PasteSpecial(-4122)
$ws.Range("F93").Value = 45922
$ws.Range("G93").Value = "O(n)"
$ws.Range("H93").Value = "O(1)"
$ws.Range("I93").Value = "2 Linked List"

# --- Row 94: Linked List Random Nodw
$ws.Range("A94").Value = 93
$ws.Range("B94").Value = "Linked List"
$ws.Range("C94").Value = "Linked List Random Nodw"
$ws.Range("D94").Value = "Medium"
$ws.Range("E94").Value = "Done"
$ws.Range("F90").Copy()
$ws.Range("F94").PasteSpecial(-4122)
$ws.Range("F94").Value = 45922
$ws.Range("G94").Value = "O(n)"
$ws.Range("H94").Value = "O(1)"
$ws.Range("I94").Value = "Using Random Class"

# --- Rows 95-101: ID-only filler cells (column A), continuing the counter
$ws.Range("A95").Value = 94
$ws.Range("A96").Value = 95
$ws.Range("A97").Value = 96
$ws.Range("A98").Value = 97
$ws.Range("A99").Value = 98
$ws.Range("A100").Value = 99
$ws.Range("A101").Value = 100

# --- Restore selection/view state to match the saved workbook
$ws.Range("A95").Select()
